$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = New-Object 'object[,]' 1,128
$values[0,0] = -0.1730910390615463
$values[0,1] = 0.0468733087182045
$values[0,2] = 0.02230261638760567
$values[0,3] = -0.1236186027526855
$values[0,4] = -0.1355118155479431
$values[0,5] = -0.004758149385452271
$values[0,6] = -0.04674947634339333
$values[0,7] = -0.08549851924180984
$values[0,8] = 0.06096422672271729
$values[0,9] = -0.1456707864999771
$values[0,10] = 0.1873953938484192
$values[0,11] = -0.05990472063422203
$values[0,12] = -0.2161078751087189
$values[0,13] = -0.07083684206008911
$values[0,14] = -0.05230434238910675
$values[0,15] = 0.1688528507947922
$values[0,16] = -0.10996925085783
$values[0,17] = -0.1554861068725586
$values[0,18] = -0.06284290552139282
$values[0,19] = 0.01328445971012115
$values[0,20] = 0.09919625520706177
$values[0,21] = 0.06926415860652924
$values[0,22] = -0.07293561846017838
$values[0,23] = 0.03993618488311768
$values[0,24] = -0.1648965626955032
$values[0,25] = -0.2873676419258118
$values[0,26] = -0.09258896857500076
$values[0,27] = 0.007343828678131104
$values[0,28] = -0.08637028932571411
$values[0,29] = -0.08937422931194305
$values[0,30] = -0.0372329093515873
$values[0,31] = 0.05746649950742722
$values[0,32] = -0.1870715022087097
$values[0,33] = 0.0144367516040802
$values[0,34] = 0.04682070761919022
$values[0,35] = 0.1388134360313416
$values[0,36] = 0.01048935949802399
$values[0,37] = -0.01035474985837936
$values[0,38] = 0.09950096905231476
$values[0,39] = 0.04611104354262352
$values[0,40] = -0.2723390758037567
$values[0,41] = 0.06730875372886658
$values[0,42] = 0.1117917224764824
$values[0,43] = 0.2708423733711243
$values[0,44] = 0.2074048519134521
$values[0,45] = -0.002255335450172424
$values[0,46] = -0.01857819408178329
$values[0,47] = -0.1428809612989426
$values[0,48] = 0.1252310872077942
$values[0,49] = -0.2293906360864639
$values[0,50] = -0.04201679304242134
$values[0,51] = 0.1312595307826996
$values[0,52] = 0.03448783606290817
$values[0,53] = 0.03538764268159866
$values[0,54] = -0.004816927015781403
$values[0,55] = -0.09243146330118179
$values[0,56] = 0.02554529905319214
$values[0,57] = 0.1159960031509399
$values[0,58] = -0.1160740330815315
$values[0,59] = 0.03537945449352264
$values[0,60] = 0.06928659975528717
$values[0,61] = -0.1114377677440643
$values[0,62] = -0.08367887139320374
$values[0,63] = -0.1209922209382057
$values[0,64] = 0.1632474362850189
$values[0,65] = 0.03589450195431709
$values[0,66] = -0.1758359968662262
$values[0,67] = -0.2392930239439011
$values[0,68] = 0.08649689704179764
$values[0,69] = -0.1891351044178009
$values[0,70] = -0.1356927454471588
$values[0,71] = 0.1131480112671852
$values[0,72] = -0.1520384252071381
$values[0,73] = -0.2002529799938202
$values[0,74] = -0.3503864407539368
$values[0,75] = 0.03609529137611389
$values[0,76] = 0.3386359214782715
$values[0,77] = 0.1422533094882965
$values[0,78] = -0.1987703293561935
$values[0,79] = 0.08287635445594788
$values[0,80] = 0.008811849169433117
$values[0,81] = 0.04027725756168365
$values[0,82] = 0.1274399161338806
$values[0,83] = 0.1737797409296036
$values[0,84] = 0.007274281233549118
$values[0,85] = 0.06602674722671509
$values[0,86] = -0.08098559081554413
$values[0,87] = 0.06098867952823639
$values[0,88] = 0.2891929149627686
$values[0,89] = -0.04164474457502365
$values[0,90] = -0.008886042982339859
$values[0,91] = 0.2377078533172607
$values[0,92] = 0.02674798667430878
$values[0,93] = 0.1161438524723053
$values[0,94] = 0.01022475212812424
$values[0,95] = 0.06723983585834503
$values[0,96] = -0.04888508468866348
$values[0,97] = -0.02901731431484222
$values[0,98] = -0.1517019271850586
$values[0,99] = 0.003211945295333862
$values[0,100] = 0.04904887825250626
$values[0,101] = -0.03636002540588379
$values[0,102] = 0.001870393753051758
$values[0,103] = 0.1718954145908356
$values[0,104] = -0.1269590258598328
$values[0,105] = 0.1904156804084778
$values[0,106] = -0.0403340682387352
$values[0,107] = 0.03302508220076561
$values[0,108] = -0.03621438145637512
$values[0,109] = -0.03502778708934784
$values[0,110] = -0.07337862998247147
$values[0,111] = 0.005858689546585083
$values[0,112] = 0.08885832130908966
$values[0,113] = -0.2003527730703354
$values[0,114] = 0.1659952402114868
$values[0,115] = 0.2087558507919312
$values[0,116] = 0.08905559778213501
$values[0,117] = 0.03161385282874107
$values[0,118] = 0.1914173662662506
$values[0,119] = -0.002605810761451721
$values[0,120] = -0.005448907613754272
$values[0,121] = 0.03450989723205566
$values[0,122] = -0.2104713171720505
$values[0,123] = -0.03083633072674274
$values[0,124] = 0.09877301007509232
$values[0,125] = -0.1010307371616364
$values[0,126] = 0.09877075254917145
$values[0,127] = -0.04287730902433395

$ws.Range("A2:DX2").Value = $values

